$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new menu item "Water" as row 8 — free (0), a Chef Special, and a Valentine Exclusive.
$ws.Range("A8").Value = "Water"
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = $true
$ws.Range("D8").Value = $true
